{"js": "// Fix a handful of typos in the requirements document:\n//   - \"meetups\" -> \"meet-ups\" (two occurrences)\n//   - \"documens\" -> \"documents\"\n//   - \"taken please\" -> \"taken place\"\n//   - \"whil each\" -> \"while each\"\n//   - \"5-star based\" -> \"5-star-based\"\n//\n// Each fix is applied with Body.search + Range.insertText(..., \"Replace\"),\n// which keeps the surrounding run formatting (font, size, color, language)\n// intact since only the matched sub-range is rewritten.\n\nasync function replaceAll(ctx, findText, replacement) {\n  const results = ctx.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await ctx.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await ctx.sync();\n}\n\nawait replaceAll(context, \"meetups\", \"meet-ups\");\nawait replaceAll(context, \"documens\", \"documents\");\nawait replaceAll(context, \"taken please\", \"taken place\");\nawait replaceAll(context, \"whil each\", \"while each\");\nawait replaceAll(context, \"5-star based\", \"5-star-based\");\n", "ps1": "# Fix a handful of typos in the requirements document:\n#   - \"meetups\" -> \"meet-ups\" (two occurrences)\n#   - \"documens\" -> \"documents\"\n#   - \"taken please\" -> \"taken place\"\n#   - \"whil each\" -> \"while each\"\n#   - \"5-star based\" -> \"5-star-based\"\n#\n# Uses Find/Replace over the whole document body (wdReplaceAll), which\n# rewrites only the matched text and leaves the run formatting (font,\n# size, color, language) of the surrounding text untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $null, $replaceText, 2) # wdReplaceAll\n}\n\nReplace-AllText \"meetups\" \"meet-ups\"\nReplace-AllText \"documens\" \"documents\"\nReplace-AllText \"taken please\" \"taken place\"\nReplace-AllText \"whil each\" \"while each\"\nReplace-AllText \"5-star based\" \"5-star-based\"\n"}
